$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(118).Insert()

$ws.Cells.Item(118, 1).Value = 8
$ws.Cells.Item(118, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 45049
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = 100112001
$ws.Cells.Item(118, 7).Value = "Berenjena"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 560
$ws.Cells.Item(118, 11).Value = 10000
$ws.Cells.Item(118, 12).Value = 11000
$ws.Cells.Item(118, 13).Value = 10500
$ws.Cells.Item(118, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(118, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(118, 16).Value = 210
$ws.Cells.Item(118, 17).Value = 50
$ws.Cells.Item(118, 18).Value = "Hortaliza"
